$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 334.5
$ws.Range("I8").Value = 334.5
$ws.Range("K8").Value = 1003.5
$ws.Range("M8").Value = -864.5
$ws.Range("H17").Value = 1974.25
$ws.Range("J17").Value = 1974.25
$ws.Range("L17").Value = 5922.75
$ws.Range("N17").Value = -6258.75
$ws.Range("H80").Value = 4830
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4830
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 14490
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -16486
$ws.Range("H83").Value = 4830
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4830
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 43470
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -53454
$ws.Range("H134").Value = 44946.8
$ws.Range("J134").Value = 44946.8
$ws.Range("L134").Value = 44946.8
$ws.Range("N134").Value = -55086.8
$ws.Range("H137").Value = 355623.53
$ws.Range("I137").Value = 1473.4286
$ws.Range("J137").Value = 2421499.2
$ws.Range("K137").Value = 4420.2858
$ws.Range("L137").Value = 7264497.600000001
$ws.Range("M137").Value = -1870.2858
$ws.Range("N137").Value = -7269597.600000001
$ws.Range("H138").Value = 2010.5217
$ws.Range("I138").Value = 1686.9166
$ws.Range("J138").Value = 2363.5454
$ws.Range("K138").Value = 5060.7498
$ws.Range("L138").Value = 7090.6362
$ws.Range("M138").Value = 79.2502000000004
$ws.Range("N138").Value = -17370.6362
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 47907.91
$ws.Range("I61").Value = 2398.8235
$ws.Range("K61").Value = 2398.8235
$ws.Range("M61").Value = -2186.8235
$ws.Range("H74").Value = 65877.94
$ws.Range("I74").Value = 85606.914
$ws.Range("K74").Value = 85606.914
$ws.Range("M74").Value = -84732.914
$ws.Range("H77").Value = 65877.94
$ws.Range("I77").Value = 85606.914
$ws.Range("K77").Value = 428034.57
$ws.Range("M77").Value = -423666.57
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H136").Value = 47907.91
$ws.Range("I136").Value = 2398.8235
$ws.Range("K136").Value = 7196.470499999999
$ws.Range("M136").Value = -4646.470499999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5742.1787
$ws.Range("I20").Value = 8718
$ws.Range("J20").Value = 1143.1818
$ws.Range("K20").Value = 8718
$ws.Range("L20").Value = 1143.1818
$ws.Range("M20").Value = -8471
$ws.Range("N20").Value = -1637.1818
$ws.Range("H132").Value = 35856.734
$ws.Range("J132").Value = 35856.734
$ws.Range("L132").Value = 35856.734
$ws.Range("N132").Value = -45976.734
$ws.Range("H135").Value = 66773.5
$ws.Range("J135").Value = 66773.5
$ws.Range("L135").Value = 66773.5
$ws.Range("N135").Value = -76913.5
$ws.Range("H140").Value = 43499
$ws.Range("J140").Value = 43499
$ws.Range("L140").Value = 43499
$ws.Range("N140").Value = -53859
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1434.8823
$ws.Range("I58").Value = 1339.5333
$ws.Range("K58").Value = 1339.5333
$ws.Range("M58").Value = -1136.5333
$ws.Range("H86").Value = 4469843.5
$ws.Range("I86").Value = 7147650
$ws.Range("K86").Value = 7147650
$ws.Range("M86").Value = -7146527
$ws.Range("H89").Value = 4469843.5
$ws.Range("I89").Value = 7147650
$ws.Range("K89").Value = 35738250
$ws.Range("M89").Value = -35732634
$ws.Range("H134").Value = 2013966.6
$ws.Range("I134").Value = 2166554.5
$ws.Range("K134").Value = 6499663.5
$ws.Range("M134").Value = -6497128.5
$ws.Range("H136").Value = 1434.8823
$ws.Range("I136").Value = 1339.5333
$ws.Range("K136").Value = 4018.5999
$ws.Range("M136").Value = -1468.5999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1077.25
$ws.Range("I5").Value = 1104.5
$ws.Range("J5").Value = 1050
$ws.Range("K5").Value = 3313.5
$ws.Range("L5").Value = 3150
$ws.Range("M5").Value = -3201.5
$ws.Range("N5").Value = -3374
$ws.Range("H132").Value = 5209
$ws.Range("J132").Value = 5454.4443
$ws.Range("L132").Value = 49089.9987
$ws.Range("N132").Value = -54149.9987
$ws.Range("H135").Value = 1077.25
$ws.Range("I135").Value = 1104.5
$ws.Range("J135").Value = 1050
$ws.Range("K135").Value = 9940.5
$ws.Range("L135").Value = 9450
$ws.Range("M135").Value = -7405.5
$ws.Range("N135").Value = -14520
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3494.5
$ws.Range("I80").Value = 3494.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3494.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2496.5
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 3494.5
$ws.Range("I83").Value = 3494.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 17472.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -12480.5
$ws.Range("N83").Value = ""
$ws.Range("H108").Value = 49240.816
$ws.Range("J108").Value = 49240.816
$ws.Range("L108").Value = 49240.816
$ws.Range("N108").Value = -56920.816
$ws.Range("H132").Value = 3923.8215
$ws.Range("I132").Value = 3146.3914
$ws.Range("K132").Value = 9439.174199999999
$ws.Range("M132").Value = -6909.174199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3040.7827
$ws.Range("I46").Value = 963.3333
$ws.Range("K46").Value = 963.3333
$ws.Range("M46").Value = -775.3333
$ws.Range("I55").Value = 622.5454999999999
$ws.Range("K55").Value = 622.5454999999999
$ws.Range("M55").Value = -449.5454999999999
$ws.Range("H136").Value = 1808.7142
$ws.Range("I136").Value = 1630.7368
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 4892.2104
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -2342.2104
$ws.Range("N136").Value = -15598.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 99999
$ws.Range("J16").Value = 99999
$ws.Range("L16").Value = 99999
$ws.Range("N16").Value = -100583
$ws.Range("H136").Value = 2059
$ws.Range("I136").Value = 2059
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6177
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3627
$ws.Range("N136").Value = ""
